$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 337
$ws.Range("F3").Value = 1170
$ws.Range("F7").Value = 175
$ws.Range("F8").Value = 669
$ws.Range("F9").Value = 1903
$ws.Range("F10").Value = 65
$ws.Range("F11").Value = 500
$ws.Range("F12").Value = 74
$ws.Range("F14").Value = 725
$ws.Range("F15").Value = 495
$ws.Range("F17").Value = 865
$ws.Range("F18").Value = 81426
$ws.Range("F19").Value = 81426
$ws.Range("F20").Value = 5
$ws.Range("F21").Value = 692
$ws.Range("F22").Value = 35040
$ws.Range("F23").Value = 35040
$ws.Range("F24").Value = 602
$ws.Range("F26").Value = 34
$ws.Range("F28").Value = 70
$ws.Range("F29").Value = 1066
$ws.Range("F30").Value = 336
$ws.Range("F32").Value = 725
$ws.Range("F33").Value = 3541
$ws.Range("F35").Value = 1277
$ws.Range("F36").Value = 5569
$ws.Range("F37").Value = 848
$ws.Range("F38").Value = 494
$ws.Range("F41").Value = 17
$ws.Range("F44").Value = 498
$ws.Range("F47").Value = 1
$ws.Range("F49").Value = 70
$ws.Range("F51").Value = 1

$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 1797
$ws.Range("F8").Value = 11
$ws.Range("F9").Value = 2014
$ws.Range("F10").Value = 45
$ws.Range("F13").Value = 429
$ws.Range("F15").Value = 81
$ws.Range("F16").Value = 80
$ws.Range("F18").Value = 561
$ws.Range("F31").Value = 506
$ws.Range("F37").Value = 24
$ws.Range("F40").Value = 43
$ws.Range("F43").Value = 76
$ws.Range("F44").Value = 846
$ws.Range("F45").Value = 313

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 739
$ws.Range("F7").Value = 228

$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 1170
$ws.Range("F10").Value = 1797
$ws.Range("F11").Value = 11
$ws.Range("F12").Value = 175
$ws.Range("F13").Value = 669
$ws.Range("F14").Value = 228
$ws.Range("F15").Value = 2014
$ws.Range("F16").Value = 1903
$ws.Range("F17").Value = 45
$ws.Range("F19").Value = 500
$ws.Range("F22").Value = 725
$ws.Range("F23").Value = 865
$ws.Range("F24").Value = 81
$ws.Range("F25").Value = 81427
$ws.Range("F26").Value = 80
$ws.Range("F27").Value = 35041
$ws.Range("F29").Value = 34
$ws.Range("F31").Value = 561
$ws.Range("F32").Value = 561
$ws.Range("F33").Value = 70
$ws.Range("F34").Value = 1066
$ws.Range("F36").Value = 336
$ws.Range("F38").Value = 1277
$ws.Range("F39").Value = 5570
$ws.Range("F40").Value = 848
$ws.Range("F45").Value = 498
$ws.Range("F49").Value = 846
$ws.Range("F50").Value = 313
$ws.Range("F52").Value = 70
